$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1545
$ws.Range("F3").Value = 3302
$ws.Range("F4").Value = 19
$ws.Range("F5").Value = 727
$ws.Range("F6").Value = 2332
$ws.Range("F7").Value = 499
$ws.Range("F8").Value = 417
$ws.Range("F9").Value = 249
$ws.Range("F10").Value = 141
$ws.Range("F11").Value = 356
$ws.Range("F12").Value = 1102
$ws.Range("F13").Value = 455
$ws.Range("F14").Value = 22
$ws.Range("F16").Value = 259
$ws.Range("F17").Value = 4770
$ws.Range("F19").Value = 1354
$ws.Range("F20").Value = 3534
$ws.Range("F22").Value = 135
$ws.Range("F23").Value = 196
$ws.Range("F24").Value = 3775
$ws.Range("F25").Value = 5152
$ws.Range("F27").Value = 983
$ws.Range("F28").Value = 567
$ws.Range("F29").Value = 3315
$ws.Range("G29").Value = 52.2
$ws.Range("F30").Value = 380
$ws.Range("F32").Value = 144
$ws.Range("F34").Value = 892
$ws.Range("F35").Value = 1208
$ws.Range("F36").Value = 20
$ws.Range("F37").Value = 22
$ws.Range("F38").Value = 1430
$ws.Range("F39").Value = 140
$ws.Range("F40").Value = 1401
$ws.Range("F41").Value = 896
$ws.Range("F42").Value = 871
$ws.Range("F43").Value = 515
$ws.Range("F44").Value = 59
$ws.Range("F45").Value = 355
$ws.Range("F46").Value = 77
$ws.Range("F47").Value = 173
$ws.Range("F49").Value = 3745

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 1022
$ws.Range("F11").Value = 17

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2323

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2323
$ws.Range("F3").Value = 1545
$ws.Range("F4").Value = 3302
$ws.Range("F5").Value = 19
$ws.Range("F6").Value = 727
$ws.Range("F7").Value = 2332
$ws.Range("F8").Value = 499
$ws.Range("F9").Value = 417
$ws.Range("F10").Value = 249
$ws.Range("F11").Value = 1022
$ws.Range("F12").Value = 141
$ws.Range("F13").Value = 356
$ws.Range("F14").Value = 1102
$ws.Range("F15").Value = 455
$ws.Range("F16").Value = 22
$ws.Range("F18").Value = 259
$ws.Range("F19").Value = 4770
$ws.Range("F20").Value = 1354
$ws.Range("F21").Value = 3534
$ws.Range("F22").Value = 3775
$ws.Range("F23").Value = 5152
$ws.Range("F24").Value = 983
$ws.Range("F25").Value = 567
$ws.Range("F26").Value = 3315
$ws.Range("G26").Value = 52.2
$ws.Range("F27").Value = 380
$ws.Range("F29").Value = 144
$ws.Range("F31").Value = 892
$ws.Range("F32").Value = 1208
$ws.Range("F33").Value = 20
$ws.Range("F34").Value = 22
$ws.Range("F35").Value = 1430
$ws.Range("F36").Value = 1401
$ws.Range("F37").Value = 896
$ws.Range("F39").Value = 515
$ws.Range("F41").Value = 59
$ws.Range("F43").Value = 355
$ws.Range("F45").Value = 77
$ws.Range("F46").Value = 173
$ws.Range("F49").Value = 3745
